$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# "Boolean" sheet: the single CSV entries for the BVTQaZ and VTQaZ travel
# modes were each split out into six separate per-vehicle-type CSV files
# (LDVs, HDVs, aircraft, rail, ships, motorbikes), inserted in place of
# the old single-file rows, keeping alphabetical order.
# -----------------------------------------------------------------------
$wsBool = $wb.Worksheets.Item("Boolean")

# Row 17 held "trans/BVTQaZ/BVTQaZ.csv" -> replace with 6 rows.
$wsBool.Rows.Item(18).Resize(5).Insert()
$wsBool.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBool.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBool.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBool.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBool.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBool.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# Row 26 (after the insert above) held "trans/VTQaZ/VTQaZ.csv" -> replace
# with 6 rows.
$wsBool.Rows.Item(27).Resize(5).Insert()
$wsBool.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBool.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBool.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBool.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBool.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBool.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# A handful of blank, formatted rows trail the data (left over from
# reviewing/expanding the list in the live workbook).
$wsBool.Range("A33:A38").Font.Size = 11

# -----------------------------------------------------------------------
# Window / tab-selection state: the workbook was left with "About" as the
# active tab, "Integer" showing A13 selected, and "Boolean" scrolled down
# with A32 selected (where the newly-inserted rows were being reviewed).
# -----------------------------------------------------------------------
$wsInt = $wb.Worksheets.Item("Integer")
$wsSub = $wb.Worksheets.Item("Subscript")
$wsAbout = $wb.Worksheets.Item("About")

$wsBool.Activate() | Out-Null
$wsBool.Range("A32").Select() | Out-Null

$wsInt.Activate() | Out-Null
$wsInt.Range("A13").Select() | Out-Null

$wsSub.Activate() | Out-Null
$wsSub.Range("A1").Select() | Out-Null

$wsAbout.Activate() | Out-Null
$wsAbout.Range("A1").Select() | Out-Null

Write-Host "edit complete"
